$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("V1").Value = "Weight"
$ws.Range("W1").Value = "Measurement Unit"
$ws.Range("V1:W1").Style = $ws.Range("U1").Style

# New data rows
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = "kg/cm"

$ws.Range("V3").Value = 6
$ws.Range("W3").Value = "lbs/in"

$ws.Range("V4").Value = 7
$ws.Range("W4").Value = "kg/cm"

# Column width for W
$ws.Range("W1").ColumnWidth = 17.85546875

# Update selection to match diff (W4 active cell)
$ws.Range("W4").Select()
